$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, pushing existing rows 21+ down by one.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new record.
$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).Value = 45133
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112031
$ws.Cells.Item(21, 7).Value = "Poroto verde"
$ws.Cells.Item(21, 8).Value = "Magnum"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 50
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 13).Value = 20000
$ws.Cells.Item(21, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(21, 15).Value = "Perú"
$ws.Cells.Item(21, 16).Value = 800
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"
